$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 45.89896166666667
$ws.Range("H2").Value = 137.696885
$ws.Range("I2").Value = 0.5698328702801992
$ws.Range("J2").Value = 0.5698328702801992
$ws.Range("M2").Value = 0.5373756666666667
$ws.Range("N2").Value = 1.612127
$ws.Range("O2").Value = 0.007472820128982582
$ws.Range("P2").Value = 0.007472820128982581
$ws.Range("Q2").Value = 24.66498512493278
$ws.Range("R2").Value = 221.984866124395
$ws.Range("S2").Value = 0.004258258543185793
$ws.Range("T2").Value = 0.004258258543185793
$ws.Range("G3").Value = 45.89896166666667
$ws.Range("H3").Value = 137.696885
$ws.Range("I3").Value = 0.5698328702801992
$ws.Range("J3").Value = 0.5698328702801992
$ws.Range("O3").Value = 0.1537223653287423
$ws.Range("P3").Value = 0.1537223653287423
$ws.Range("Q3").Value = 507.3800504708695
$ws.Range("R3").Value = 4566.420454237826
$ws.Range("S3").Value = 0.08759605666153863
$ws.Range("T3").Value = 0.08759605666153862
$ws.Range("G4").Value = 45.89896166666667
$ws.Range("H4").Value = 137.696885
$ws.Range("I4").Value = 0.5698328702801992
$ws.Range("J4").Value = 0.5698328702801992
$ws.Range("M4").Value = 30.561198
$ws.Range("N4").Value = 91.683594
$ws.Range("O4").Value = 0.4249882340167162
$ws.Range("P4").Value = 0.4249882340167161
$ws.Range("Q4").Value = 1402.72725548941
$ws.Range("R4").Value = 12624.54529940469
$ws.Range("S4").Value = 0.2421722652250584
$ws.Range("T4").Value = 0.2421722652250584
$ws.Range("G5").Value = 45.89896166666667
$ws.Range("H5").Value = 137.696885
$ws.Range("I5").Value = 0.5698328702801992
$ws.Range("J5").Value = 0.5698328702801992
$ws.Range("M5").Value = 29.75783666666667
$ws.Range("N5").Value = 89.27351
$ws.Range("O5").Value = 0.4138165805255589
$ws.Range("P5").Value = 0.4138165805255589
$ws.Range("Q5").Value = 1365.853804446261
$ws.Range("R5").Value = 12292.68424001635
$ws.Range("S5").Value = 0.2358062898504164
$ws.Range("T5").Value = 0.2358062898504164
$ws.Range("I6").Value = 0.04736372570041834
$ws.Range("J6").Value = 0.04736372570041834
$ws.Range("M6").Value = 0.5373756666666667
$ws.Range("N6").Value = 1.612127
$ws.Range("O6").Value = 0.007472820128982582
$ws.Range("P6").Value = 0.007472820128982581
$ws.Range("Q6").Value = 2.050119694372444
$ws.Range("R6").Value = 18.451077249352
$ws.Range("S6").Value = 0.0003539406027976958
$ws.Range("T6").Value = 0.0003539406027976958
$ws.Range("I7").Value = 0.04736372570041834
$ws.Range("J7").Value = 0.04736372570041834
$ws.Range("O7").Value = 0.1537223653287423
$ws.Range("P7").Value = 0.1537223653287423
$ws.Range("S7").Value = 0.007280863945450051
$ws.Range("T7").Value = 0.00728086394545005
$ws.Range("I8").Value = 0.04736372570041834
$ws.Range("J8").Value = 0.04736372570041834
$ws.Range("M8").Value = 30.561198
$ws.Range("N8").Value = 91.683594
$ws.Range("O8").Value = 0.4249882340167162
$ws.Range("P8").Value = 0.4249882340167161
$ws.Range("Q8").Value = 116.592763293616
$ws.Range("R8").Value = 1049.334869642544
$ws.Range("S8").Value = 0.02012902614187295
$ws.Range("T8").Value = 0.02012902614187295
$ws.Range("I9").Value = 0.04736372570041834
$ws.Range("J9").Value = 0.04736372570041834
$ws.Range("M9").Value = 29.75783666666667
$ws.Range("N9").Value = 89.27351
$ws.Range("O9").Value = 0.4138165805255589
$ws.Range("P9").Value = 0.4138165805255589
$ws.Range("Q9").Value = 113.5278926764178
$ws.Range("R9").Value = 1021.75103408776
$ws.Range("S9").Value = 0.01959989501029765
$ws.Range("T9").Value = 0.01959989501029765
$ws.Range("G10").Value = 4.651706333333334
$ws.Range("H10").Value = 13.955119
$ws.Range("I10").Value = 0.05775065655894644
$ws.Range("J10").Value = 0.05775065655894644
$ws.Range("M10").Value = 0.5373756666666667
$ws.Range("N10").Value = 1.612127
$ws.Range("O10").Value = 0.007472820128982582
$ws.Range("P10").Value = 0.007472820128982581
$ws.Range("Q10").Value = 2.499713792012556
$ws.Range("R10").Value = 22.497424128113
$ws.Range("S10").Value = 0.0004315602687956549
$ws.Range("T10").Value = 0.0004315602687956549
$ws.Range("G11").Value = 4.651706333333334
$ws.Range("H11").Value = 13.955119
$ws.Range("I11").Value = 0.05775065655894644
$ws.Range("J11").Value = 0.05775065655894644
$ws.Range("O11").Value = 0.1537223653287423
$ws.Range("P11").Value = 0.1537223653287423
$ws.Range("Q11").Value = 51.42127203928391
$ws.Range("R11").Value = 462.7914483535551
$ws.Range("S11").Value = 0.008877567525529095
$ws.Range("T11").Value = 0.008877567525529094
$ws.Range("G12").Value = 4.651706333333334
$ws.Range("H12").Value = 13.955119
$ws.Range("I12").Value = 0.05775065655894644
$ws.Range("J12").Value = 0.05775065655894644
$ws.Range("M12").Value = 30.561198
$ws.Range("N12").Value = 91.683594
$ws.Range("O12").Value = 0.4249882340167162
$ws.Range("P12").Value = 0.4249882340167161
$ws.Range("Q12").Value = 142.161718290854
$ws.Range("R12").Value = 1279.455464617686
$ws.Range("S12").Value = 0.02454334954429254
$ws.Range("T12").Value = 0.02454334954429253
$ws.Range("G13").Value = 4.651706333333334
$ws.Range("H13").Value = 13.955119
$ws.Range("I13").Value = 0.05775065655894644
$ws.Range("J13").Value = 0.05775065655894644
$ws.Range("M13").Value = 29.75783666666667
$ws.Range("N13").Value = 89.27351
$ws.Range("O13").Value = 0.4138165805255589
$ws.Range("P13").Value = 0.4138165805255589
$ws.Range("Q13").Value = 138.4247172886323
$ws.Range("R13").Value = 1245.82245559769
$ws.Range("S13").Value = 0.02389817922032916
$ws.Range("T13").Value = 0.02389817922032916
$ws.Range("G14").Value = 26.182385
$ws.Range("H14").Value = 78.547155
$ws.Range("I14").Value = 0.325052747460436
$ws.Range("J14").Value = 0.325052747460436
$ws.Range("M14").Value = 0.5373756666666667
$ws.Range("N14").Value = 1.612127
$ws.Range("O14").Value = 0.007472820128982582
$ws.Range("P14").Value = 0.007472820128982581
$ws.Range("Q14").Value = 14.06977659429833
$ws.Range("R14").Value = 126.627989348685
$ws.Range("S14").Value = 0.002429060714203438
$ws.Range("T14").Value = 0.002429060714203438
$ws.Range("G15").Value = 26.182385
$ws.Range("H15").Value = 78.547155
$ws.Range("I15").Value = 0.325052747460436
$ws.Range("J15").Value = 0.325052747460436
$ws.Range("O15").Value = 0.1537223653287423
$ws.Range("P15").Value = 0.1537223653287423
$ws.Range("Q15").Value = 289.4274584951084
$ws.Range("R15").Value = 2604.847126455975
$ws.Range("S15").Value = 0.04996787719622457
$ws.Range("T15").Value = 0.04996787719622456
$ws.Range("G16").Value = 26.182385
$ws.Range("H16").Value = 78.547155
$ws.Range("I16").Value = 0.325052747460436
$ws.Range("J16").Value = 0.325052747460436
$ws.Range("M16").Value = 30.561198
$ws.Range("N16").Value = 91.683594
$ws.Range("O16").Value = 0.4249882340167162
$ws.Range("P16").Value = 0.4249882340167161
$ws.Range("Q16").Value = 800.16505209723
$ws.Range("R16").Value = 7201.48546887507
$ws.Range("S16").Value = 0.1381435931054923
$ws.Range("T16").Value = 0.1381435931054923
$ws.Range("G17").Value = 26.182385
$ws.Range("H17").Value = 78.547155
$ws.Range("I17").Value = 0.325052747460436
$ws.Range("J17").Value = 0.325052747460436
$ws.Range("M17").Value = 29.75783666666667
$ws.Range("N17").Value = 89.27351
$ws.Range("O17").Value = 0.4138165805255589
$ws.Range("P17").Value = 0.4138165805255589
$ws.Range("Q17").Value = 779.1311363737833
$ws.Range("R17").Value = 7012.180227364051
$ws.Range("S17").Value = 0.1345122164445157
$ws.Range("T17").Value = 0.1345122164445157
